$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.155.38"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "1.568.89"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  +0.68%  "
$ws.Range("D5").Value = "'211.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.13%  "
$ws.Range("D6").Value = "'0.492"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("E7").Value = "  +0.73%  "
$ws.Range("D8").Value = "'21.99"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("D10").Value = "'0.0599"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("D11").Value = "'0.0866"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("D12").Value = "1.790.18"
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D13").Value = "1.545.02"
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").Value = "'0.519"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").Value = "27.107.93"
$ws.Range("D17").Value = "'62.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").Value = "0.0₃0703"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").Value = "'215.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.70%  "
$ws.Range("D20").Value = "'7.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("D22").Value = "'4.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.29%  "
$ws.Range("D23").Value = "'9.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("D25").Value = "'154.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("D26").Value = "'6.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").Value = "'15.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("E28").Value = "  +1.85%  "
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("E30").Value = "  +3.98%  "
$ws.Range("D31").Value = "'0.0472"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("D32").Value = "'3.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").Value = "'3.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.49%  "
$ws.Range("D34").Value = "1.455.49"
$ws.Range("E34").Value = "  +2.30%  "
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("D39").Value = "'0.533"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.48%  "
$ws.Range("D40").Value = "'5.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.39%  "
$ws.Range("D41").Value = "'0.809"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("E42").Value = "  +0.80%  "
$ws.Range("E43").Value = "  +0.78%  "
$ws.Range("E44").Value = "  -2.16%  "
$ws.Range("D45").Value = "'64.62"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("D47").Value = "1.701.46"
$ws.Range("E47").Value = "  +0.52%  "
$ws.Range("D48").Value = "'86.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.63%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0103"
$ws.Range("E49").Value = "  +2.64%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0519"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("D51").Value = "'0.0959"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.02%  "
